# Data read from excel sheet in post request method implemented.
#
# "Sheet1" holds the name/job values that the POST-request test cases pull
# from the spreadsheet. Row 3 now supplies a real user name ("Urmila") with
# the job left as "tester", and row 4 is switched to exercise the numeric
# name (123) with job "tester" - replacing the old placeholder rows
# (123/456 and the special-character "$%*"/"@&^" pair).
#
# Sheet1 becomes the active/selected tab (it was "Test Cases" before).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Test Cases"
$ws2 = $wb.Worksheets.Item(2)   # "Sheet1"

# Row 3: was (123, 456) -> now ("Urmila", "tester")
$ws2.Range("A3").Value = "Urmila"
$ws2.Range("B3").Value = "tester"

# Row 4: was ("$%*", "@&^") -> now (123, "tester")
$ws2.Range("A4").Value = 123
$ws2.Range("B4").Value = "tester"

# Make "Sheet1" the active tab/window, with D5 selected (was E5).
$ws2.Activate()
$ws2.Range("D5").Select()
